# Added a customer PUT method for billing address
# Insert a new "id" column at the front of the Billing sheet, then
# select it and make Billing the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Billing")

# Insert a new column before column A, shifting existing data right.
$ws.Columns("A:A").Insert()

# Populate the new column.
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 155

# Update the selection and make Billing the active sheet/tab.
$ws.Activate()
$ws.Range("B2").Select()
